$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

$ws.Cells.Item(2, 1).Value = "Lenovo Ideapad S145 AMD A6-9225 15.6 inch HD Thin and Light Laptop (4GB/1TB/Windows 10/Grey/1.85Kg), 81N30063IN"
$ws.Cells.Item(3, 1).Value = "Lenovo V145-AMD-A6 15.6 inch HD Thin and Light Laptop (4GB RAM/ 500GB HDD/ Windows 10 Home with Lifetime Validity/ Black/ 2.1 kg), 81MT004BIH"
$ws.Cells.Item(4, 1).Value = "Lenovo Ideapad S145 7th Gen Intel Core i3 15.6-inch FHD Thin and Light Laptop (4GB/1TB HDD/Windows 10/Textured Black/1.85Kg), 81VD002YIN"
$ws.Cells.Item(5, 1).Value = "Lenovo Ideapad S145 7th Gen Intel Core i3 15.6 inch FHD Thin and Light Laptop (4GB/1TB/Windows 10/Grey/1.85Kg), 81VD0008IN"
$ws.Cells.Item(6, 1).Value = "Lenovo Ideapad S145 Intel Core I3 8th Gen 15.6-inch FHD Thin and Light Laptop ( 8GB RAM / 1TB HDD / DOS / Black / 1.85 Kg),81MV0094IN"
$ws.Cells.Item(7, 1).Value = "Lenovo Ideapad S145 Intel Core I3 8th Gen 15.6-inch Thin and Light FHD Laptop ( 4GB RAM / 1TB HDD / Windows 10 Home / Grey / 1.85Kg ), 81MV0091IN"
$ws.Cells.Item(8, 1).Value = "Lenovo Ideapad 330 81DE0363IN 15.6-inch FHD Laptop (8th Gen I5-8250U/8GB/1TB HDD/Windows 10/Integrated Graphics), Platinum Grey"
$ws.Cells.Item(9, 1).Value = "Lenovo IdeaPad S145 81W800C3IN 15.6-inch FHD  Thin and Light Laptop (10th Gen CORE I3-1005G1/4GB/256GB SSD/Windows 10/Microsoft Office/Integrated Graphics), Grey"
$ws.Cells.Item(10, 1).Value = "Lenovo Ideapad S145 8th Gen Intel Core I5 15.6 inch FHD Thin and Light Laptop (8 GB RAM/ 1 TB HDD/ Windows 10/ Glossy Black / 1.85 Kg), 81MV0098IN"
$ws.Cells.Item(11, 1).Value = "Lenovo Legion Y540 9th Gen Intel Core i5 15.6 inch FHD Gaming Laptop -Lenovo 2TB External Hard Drive"
$ws.Cells.Item(12, 1).Value = "Lenovo IdeaPad S340 81VW00CVIN 15.6-inch FHD IPS Thin and Light Laptop (10th Gen CORE I5-1035G4/8GB/512GB SSD/Windows 10/Microsoft Office/Integrated Graphics), Platinum Grey"
$ws.Cells.Item(13, 1).Value = "Lenovo Ideapad 330 AMD A6-9225 Processor 15.6-inch HD Laptop (4GB/1TB HDD/DOS/Onyx Black/2.2Kg), 81D60079IN"
$ws.Cells.Item(14, 1).Value = "Lenovo IdeaPad S145 8th Gen Intel Core i5 15.6-inch FHD Thin and Light Laptop (8GB/1TB/DOS/Textured Black/1.85Kg), 81MV0166IN"
$ws.Cells.Item(15, 1).Value = "Lenovo Ideapad S145 81N300F2IN 15.6-inch HD Thin and Light Laptop (7th Gen A6-9225/4GB/1TB HDD/DOS/Integrated Graphics), Grey"
$ws.Cells.Item(16, 1).Value = "Lenovo Ideapad 330 Intel Core i5 8th Gen 15.6-inch Full HD Laptop (8GB DDR4/1TB HDD/Windows 10 Home/Platinum Grey/ 2.2kg), 81DE008PIN"
$ws.Cells.Item(17, 1).Value = "Lenovo IdeaPad S145 AMD A6 -9225 15.6-inch HD Thin and Light Laptop (4GB/1TB/Windows 10/MS Office 2019/Textured Black/1.85Kg), 81N300B7IN"
